$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C26/D26: cost values for M3x24 standoffs (match currency style used by sibling rows)
$ws.Range("C27:D27").Copy()
$ws.Range("C26:D26").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 9

# E26: replace "CANNOT FIND " text with a hyperlink to the new URL
$url = "https://www.amazon.com/jing-Standoff-Quadcopter-Computer-Circuit/dp/B0975SLD2K/ref=sr_1_1_sspa?keywords=M3+x+25+standoff&qid=1643141993&s=industrial&sr=1-1-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUExQkg1MUNCTk9GWENHJmVuY3J5cHRlZElkPUEwNTU3MzkxMUpJUDBXM0E2SE9YOSZlbmNyeXB0ZWRBZElkPUEwOTkwNzE1MkoxVTFFUkxGVUFHOSZ3aWRnZXROYW1lPXNwX2F0ZiZhY3Rpb249Y2xpY2tSZWRpcmVjdCZkb05vdExvZ0NsaWNrPXRydWU="

$ws.Hyperlinks.Add($ws.Range("E26"), $url, [Type]::Missing, [Type]::Missing, $url)
$ws.Range("E26").Style = "Hyperlink"

# Row 35: Total row
$ws.Range("C35").Value = "Total"
$ws.Range("D35").Formula = "=SUM(D4:D31)"
